$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"

$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108007
$ws.Cells.Item($row, 10).Value = "Coco"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 20000
$ws.Cells.Item($row, 15).Value = 21000
$ws.Cells.Item($row, 16).Value = 20500
$ws.Cells.Item($row, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item($row, 18).Value = "Perú"
$ws.Cells.Item($row, 19).Value = 1025
$ws.Cells.Item($row, 20).Value = 20
